$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = "kyu"
$ws.Range("B4").Value = "why"
$ws.Range("C4").Value = "'1234567890"

# Row 5
$ws.Range("A5").Value = "kyu"
$ws.Range("B5").Value = "why"
$ws.Range("C5").Value = "'12345678566"
